# Generate Report for Handoff
# Refresh the latest-handoff timestamps now that the de-de handoff has
# completed: the Overview sheet's "Latest Handoff Date" and the de-de
# sheet's "Latest Handoff Datetime" both move to the new handoff time,
# and the zh-cn sheet's "Latest Handoff Datetime" is refreshed too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D2").Value = "2016-03-13 09:03:15"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 09:03:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 09:03:15"
